$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy styles for column A (bold/border/center) and column E (date format) from the last existing row (196) down through the new rows
$ws.Cells.Item(196, 1).Copy($ws.Range($ws.Cells.Item(197, 1), $ws.Cells.Item(206, 1)))
$ws.Cells.Item(196, 5).Copy($ws.Range($ws.Cells.Item(197, 5), $ws.Cells.Item(206, 5)))

# Row 197
$ws.Cells.Item(197, 1).Value = 195
$ws.Cells.Item(197, 2).Value = 7863104
$ws.Cells.Item(197, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(197, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(197, 5).Value = 45350.60416666666
$ws.Cells.Item(197, 6).Value = "Stal Mielec"
$ws.Cells.Item(197, 7).Value = "LKS Lodz"
$ws.Cells.Item(197, 8).Value = 1
$ws.Cells.Item(197, 9).Value = 0
$ws.Cells.Item(197, 10).Value = "H"
$ws.Cells.Item(197, 11).Value = 1.727
$ws.Cells.Item(197, 12).Value = 3.75
$ws.Cells.Item(197, 13).Value = 3.75
$ws.Cells.Item(197, 14).Value = 2
$ws.Cells.Item(197, 15).Value = 3.3
$ws.Cells.Item(197, 16).Value = 3.5
$ws.Cells.Item(197, 17).Value = -0.5
$ws.Cells.Item(197, 18).Value = 2.05
$ws.Cells.Item(197, 19).Value = 1.8
$ws.Cells.Item(197, 20).Value = 2.25
$ws.Cells.Item(197, 21).Value = 1.925
$ws.Cells.Item(197, 22).Value = 1.925
$ws.Cells.Item(197, 23).Value = 1
$ws.Cells.Item(197, 24).Value = -1
$ws.Cells.Item(197, 25).Value = -1
$ws.Cells.Item(197, 26).Value = 1.05
$ws.Cells.Item(197, 27).Value = -1
$ws.Cells.Item(197, 28).Value = -1
$ws.Cells.Item(197, 29).Value = 0.925

# Row 198
$ws.Cells.Item(198, 1).Value = 196
$ws.Cells.Item(198, 2).Value = 6775554
$ws.Cells.Item(198, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(198, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(198, 5).Value = 45352.58333333334
$ws.Cells.Item(198, 6).Value = "Cracovia Krakow"
$ws.Cells.Item(198, 7).Value = "Warta Poznan"
$ws.Cells.Item(198, 11).Value = 1.8
$ws.Cells.Item(198, 12).Value = 3.3
$ws.Cells.Item(198, 13).Value = 4.333
$ws.Cells.Item(198, 14).Value = 1.75
$ws.Cells.Item(198, 15).Value = 3.3
$ws.Cells.Item(198, 16).Value = 4.5
$ws.Cells.Item(198, 17).Value = -0.75
$ws.Cells.Item(198, 18).Value = 2.05
$ws.Cells.Item(198, 19).Value = 1.8
$ws.Cells.Item(198, 20).Value = 2.25
$ws.Cells.Item(198, 21).Value = 2.05
$ws.Cells.Item(198, 22).Value = 1.8
$ws.Cells.Item(198, 23).Value = 0
$ws.Cells.Item(198, 24).Value = 0
$ws.Cells.Item(198, 25).Value = 0
$ws.Cells.Item(198, 26).Value = 0
$ws.Cells.Item(198, 27).Value = 0

# Row 199
$ws.Cells.Item(199, 1).Value = 197
$ws.Cells.Item(199, 2).Value = 6774463
$ws.Cells.Item(199, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(199, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(199, 5).Value = 45352.6875
$ws.Cells.Item(199, 6).Value = "Ruch Chorzow"
$ws.Cells.Item(199, 7).Value = "Piast Gliwice"
$ws.Cells.Item(199, 11).Value = 3.6
$ws.Cells.Item(199, 12).Value = 3
$ws.Cells.Item(199, 13).Value = 2.05
$ws.Cells.Item(199, 14).Value = 3.6
$ws.Cells.Item(199, 15).Value = 3
$ws.Cells.Item(199, 16).Value = 2.05
$ws.Cells.Item(199, 17).Value = 0.25
$ws.Cells.Item(199, 18).Value = 2
$ws.Cells.Item(199, 19).Value = 1.85
$ws.Cells.Item(199, 20).Value = 2
$ws.Cells.Item(199, 21).Value = 1.85
$ws.Cells.Item(199, 22).Value = 2
$ws.Cells.Item(199, 23).Value = 0
$ws.Cells.Item(199, 24).Value = 0
$ws.Cells.Item(199, 25).Value = 0
$ws.Cells.Item(199, 26).Value = 0
$ws.Cells.Item(199, 27).Value = 0

# Row 200
$ws.Cells.Item(200, 1).Value = 198
$ws.Cells.Item(200, 2).Value = 6775555
$ws.Cells.Item(200, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(200, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(200, 5).Value = 45353.45833333334
$ws.Cells.Item(200, 6).Value = "Gornik Zabrze"
$ws.Cells.Item(200, 7).Value = "Jagiellonia Bialystok"
$ws.Cells.Item(200, 11).Value = 2.55
$ws.Cells.Item(200, 12).Value = 3.3
$ws.Cells.Item(200, 13).Value = 2.45
$ws.Cells.Item(200, 14).Value = 2.4
$ws.Cells.Item(200, 15).Value = 3.3
$ws.Cells.Item(200, 16).Value = 2.625
$ws.Cells.Item(200, 17).Value = 0
$ws.Cells.Item(200, 18).Value = 1.85
$ws.Cells.Item(200, 19).Value = 2
$ws.Cells.Item(200, 20).Value = 2.5
$ws.Cells.Item(200, 21).Value = 1.825
$ws.Cells.Item(200, 22).Value = 2.025
$ws.Cells.Item(200, 23).Value = 0
$ws.Cells.Item(200, 24).Value = 0
$ws.Cells.Item(200, 25).Value = 0
$ws.Cells.Item(200, 26).Value = 0
$ws.Cells.Item(200, 27).Value = 0

# Row 201
$ws.Cells.Item(201, 1).Value = 199
$ws.Cells.Item(201, 2).Value = 6775560
$ws.Cells.Item(201, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(201, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(201, 5).Value = 45353.5625
$ws.Cells.Item(201, 6).Value = "Slask Wroclaw"
$ws.Cells.Item(201, 7).Value = "Widzew Lodz"
$ws.Cells.Item(201, 11).Value = 2.1
$ws.Cells.Item(201, 12).Value = 3.2
$ws.Cells.Item(201, 13).Value = 3.25
$ws.Cells.Item(201, 14).Value = 2.05
$ws.Cells.Item(201, 15).Value = 3.2
$ws.Cells.Item(201, 16).Value = 3.4
$ws.Cells.Item(201, 17).Value = -0.25
$ws.Cells.Item(201, 18).Value = 1.8
$ws.Cells.Item(201, 19).Value = 2.05
$ws.Cells.Item(201, 20).Value = 2.25
$ws.Cells.Item(201, 21).Value = 1.925
$ws.Cells.Item(201, 22).Value = 1.925
$ws.Cells.Item(201, 23).Value = 0
$ws.Cells.Item(201, 24).Value = 0
$ws.Cells.Item(201, 25).Value = 0
$ws.Cells.Item(201, 26).Value = 0
$ws.Cells.Item(201, 27).Value = 0

# Row 202
$ws.Cells.Item(202, 1).Value = 200
$ws.Cells.Item(202, 2).Value = 6775557
$ws.Cells.Item(202, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(202, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(202, 5).Value = 45353.66666666666
$ws.Cells.Item(202, 6).Value = "Legia Warsaw"
$ws.Cells.Item(202, 7).Value = "Pogon Szczecin"
$ws.Cells.Item(202, 11).Value = 2
$ws.Cells.Item(202, 12).Value = 3.5
$ws.Cells.Item(202, 13).Value = 3.3
$ws.Cells.Item(202, 14).Value = 2
$ws.Cells.Item(202, 15).Value = 3.5
$ws.Cells.Item(202, 16).Value = 3.3
$ws.Cells.Item(202, 17).Value = -0.5
$ws.Cells.Item(202, 18).Value = 2.05
$ws.Cells.Item(202, 19).Value = 1.8
$ws.Cells.Item(202, 20).Value = 2.75
$ws.Cells.Item(202, 21).Value = 2
$ws.Cells.Item(202, 22).Value = 1.85
$ws.Cells.Item(202, 23).Value = 0
$ws.Cells.Item(202, 24).Value = 0
$ws.Cells.Item(202, 25).Value = 0
$ws.Cells.Item(202, 26).Value = 0
$ws.Cells.Item(202, 27).Value = 0

# Row 203
$ws.Cells.Item(203, 1).Value = 201
$ws.Cells.Item(203, 2).Value = 6774880
$ws.Cells.Item(203, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(203, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(203, 5).Value = 45354.35416666666
$ws.Cells.Item(203, 6).Value = "LKS Lodz"
$ws.Cells.Item(203, 7).Value = "MKS Puszcza Niepolomice"
$ws.Cells.Item(203, 11).Value = 2.45
$ws.Cells.Item(203, 12).Value = 3.3
$ws.Cells.Item(203, 13).Value = 2.55
$ws.Cells.Item(203, 14).Value = 2.375
$ws.Cells.Item(203, 15).Value = 3.3
$ws.Cells.Item(203, 16).Value = 2.625
$ws.Cells.Item(203, 17).Value = 0
$ws.Cells.Item(203, 18).Value = 1.775
$ws.Cells.Item(203, 19).Value = 2.1
$ws.Cells.Item(203, 20).Value = 2.5
$ws.Cells.Item(203, 21).Value = 1.975
$ws.Cells.Item(203, 22).Value = 1.875
$ws.Cells.Item(203, 23).Value = 0
$ws.Cells.Item(203, 24).Value = 0
$ws.Cells.Item(203, 25).Value = 0
$ws.Cells.Item(203, 26).Value = 0
$ws.Cells.Item(203, 27).Value = 0

# Row 204
$ws.Cells.Item(204, 1).Value = 202
$ws.Cells.Item(204, 2).Value = 6775556
$ws.Cells.Item(204, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(204, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(204, 5).Value = 45354.45833333334
$ws.Cells.Item(204, 6).Value = "Zaglebie Lubin"
$ws.Cells.Item(204, 7).Value = "Korona Kielce"
$ws.Cells.Item(204, 11).Value = 2.05
$ws.Cells.Item(204, 12).Value = 3.2
$ws.Cells.Item(204, 13).Value = 3.2
$ws.Cells.Item(204, 14).Value = 2.05
$ws.Cells.Item(204, 15).Value = 3.2
$ws.Cells.Item(204, 16).Value = 3.2
$ws.Cells.Item(204, 17).Value = -0.25
$ws.Cells.Item(204, 18).Value = 1.85
$ws.Cells.Item(204, 19).Value = 2
$ws.Cells.Item(204, 20).Value = 2.5
$ws.Cells.Item(204, 21).Value = 2.05
$ws.Cells.Item(204, 22).Value = 1.8
$ws.Cells.Item(204, 23).Value = 0
$ws.Cells.Item(204, 24).Value = 0
$ws.Cells.Item(204, 25).Value = 0
$ws.Cells.Item(204, 26).Value = 0
$ws.Cells.Item(204, 27).Value = 0

# Row 205
$ws.Cells.Item(205, 1).Value = 203
$ws.Cells.Item(205, 2).Value = 6775559
$ws.Cells.Item(205, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(205, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(205, 5).Value = 45354.5625
$ws.Cells.Item(205, 6).Value = "Rakow Czestochowa"
$ws.Cells.Item(205, 7).Value = "Lech Poznan"
$ws.Cells.Item(205, 11).Value = 1.95
$ws.Cells.Item(205, 12).Value = 3.25
$ws.Cells.Item(205, 13).Value = 3.8
$ws.Cells.Item(205, 14).Value = 1.95
$ws.Cells.Item(205, 15).Value = 3.25
$ws.Cells.Item(205, 16).Value = 3.8
$ws.Cells.Item(205, 17).Value = -0.5
$ws.Cells.Item(205, 18).Value = 2
$ws.Cells.Item(205, 19).Value = 1.85
$ws.Cells.Item(205, 20).Value = 2.25
$ws.Cells.Item(205, 21).Value = 1.9
$ws.Cells.Item(205, 22).Value = 1.95
$ws.Cells.Item(205, 23).Value = 0
$ws.Cells.Item(205, 24).Value = 0
$ws.Cells.Item(205, 25).Value = 0
$ws.Cells.Item(205, 26).Value = 0
$ws.Cells.Item(205, 27).Value = 0

# Row 206
$ws.Cells.Item(206, 1).Value = 204
$ws.Cells.Item(206, 2).Value = 6775558
$ws.Cells.Item(206, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(206, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(206, 5).Value = 45355.625
$ws.Cells.Item(206, 6).Value = "Radomiak Radom"
$ws.Cells.Item(206, 7).Value = "Stal Mielec"
$ws.Cells.Item(206, 11).Value = 2.05
$ws.Cells.Item(206, 12).Value = 3.2
$ws.Cells.Item(206, 13).Value = 3.5
$ws.Cells.Item(206, 14).Value = 1.85
$ws.Cells.Item(206, 15).Value = 3.3
$ws.Cells.Item(206, 16).Value = 4
$ws.Cells.Item(206, 17).Value = -0.5
$ws.Cells.Item(206, 18).Value = 1.9
$ws.Cells.Item(206, 19).Value = 1.95
$ws.Cells.Item(206, 20).Value = 2.25
$ws.Cells.Item(206, 21).Value = 1.9
$ws.Cells.Item(206, 22).Value = 1.95
$ws.Cells.Item(206, 23).Value = 0
$ws.Cells.Item(206, 24).Value = 0
$ws.Cells.Item(206, 25).Value = 0
$ws.Cells.Item(206, 26).Value = 0
$ws.Cells.Item(206, 27).Value = 0

Write-Output "done"